$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 75 ("Choclero" / "Segunda",
# Región del Maule, fecha 2022-01-05). All the previously existing data rows
# 75..131 shift down by one (to 76..132); nothing else about them changes.
$ws.Rows.Item(75).Insert()

$ws.Cells.Item(75, 1).Value  = 7
$ws.Cells.Item(75, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value  = "Ñuble"
$ws.Cells.Item(75, 4).Value  = 44566
$ws.Cells.Item(75, 5).Value  = 16
$ws.Cells.Item(75, 6).Value  = 100112024
$ws.Cells.Item(75, 7).Value  = "Choclo"
$ws.Cells.Item(75, 8).Value  = "Choclero"
$ws.Cells.Item(75, 9).Value  = "Segunda"
$ws.Cells.Item(75, 10).Value = 2000
$ws.Cells.Item(75, 11).Value = 200
$ws.Cells.Item(75, 12).Value = 250
$ws.Cells.Item(75, 13).Value = 225
$ws.Cells.Item(75, 14).Value = "$/unidad"
$ws.Cells.Item(75, 15).Value = "Región del Maule"
$ws.Cells.Item(75, 16).Value = 225
$ws.Cells.Item(75, 17).Value = 1
$ws.Cells.Item(75, 18).Value = "Hortaliza"
